$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set "Army Purple" as the K's Color for every part whose Ori Color is "purple"
$rows = @(11, 12, 14, 20, 24, 25)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Army Purple"
}

# Update the selected cell to match the saved view state
$ws.Range("C23").Select()
